$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029619739271163
$ws.Range("D2").Value = 1.052000944284393
$ws.Range("E2").Value = 1.02942828623817
$ws.Range("F2").Value = 1.055843262762643
$ws.Range("I2").Value = 1.039968330270134
$ws.Range("J2").Value = 1.034765509363356
$ws.Range("K2").Value = 1.054750886844637
$ws.Range("L2").Value = 1.032242225918404
$ws.Range("M2").Value = 1.05858260710631
$ws.Range("N2").Value = 1.036234995075915
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030503058064487
$ws.Range("D3").Value = 1.05269579014621
$ws.Range("E3").Value = 1.03017709022274
$ws.Range("F3").Value = 1.056719449862806
$ws.Range("I3").Value = 1.040173867365555
$ws.Range("J3").Value = 1.035290275148491
$ws.Range("K3").Value = 1.055258794236249
$ws.Range("L3").Value = 1.032799333008386
$ws.Range("M3").Value = 1.059272155921678
$ws.Range("N3").Value = 1.036760506088657
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031075254845083
$ws.Range("D4").Value = 1.053145797839562
$ws.Range("E4").Value = 1.030662535828556
$ws.Range("F4").Value = 1.05728730362465
$ws.Range("I4").Value = 1.040305940727859
$ws.Range("J4").Value = 1.035629830602661
$ws.Range("K4").Value = 1.0555871620948
$ws.Range("L4").Value = 1.033160073788991
$ws.Range("M4").Value = 1.05971860065778
$ws.Range("N4").Value = 1.037100543750519
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031315955955575
$ws.Range("D5").Value = 1.053335074103428
$ws.Range("E5").Value = 1.030866835635013
$ws.Range("F5").Value = 1.057526243474572
$ws.Range("I5").Value = 1.040361242745122
$ws.Range("J5").Value = 1.035772577842886
$ws.Range("K5").Value = 1.055725139035968
$ws.Range("L5").Value = 1.033311789087911
$ws.Range("M5").Value = 1.059906346821994
$ws.Range("N5").Value = 1.037243493708198
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031356379463736
$ws.Range("D6").Value = 1.053366859824999
$ws.Range("E6").Value = 1.030901151235722
$ws.Range("F6").Value = 1.057566375029862
$ws.Range("I6").Value = 1.040370515181153
$ws.Range("J6").Value = 1.035796545589973
$ws.Range("K6").Value = 1.055748301891553
$ws.Range("L6").Value = 1.033337266229505
$ws.Range("M6").Value = 1.05993787374842
$ws.Range("N6").Value = 1.037267495492235
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031078470519017
$ws.Range("D7").Value = 1.053148326594936
$ws.Range("E7").Value = 1.030665264836152
$ws.Range("F7").Value = 1.057290495511549
$ws.Range("I7").Value = 1.040306680547804
$ws.Range("J7").Value = 1.035631738007044
$ws.Range("K7").Value = 1.055589006021889
$ws.Range("L7").Value = 1.033162100782515
$ws.Range("M7").Value = 1.0597211090943
$ws.Range("N7").Value = 1.037102453863634
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029918129641046
$ws.Range("D8").Value = 1.052235687650947
$ws.Range("E8").Value = 1.029681156716595
$ws.Range("F8").Value = 1.056139186398724
$ws.Range("I8").Value = 1.04003798317072
$ws.Range("J8").Value = 1.034942856378786
$ws.Range("K8").Value = 1.054922594142109
$ws.Range("L8").Value = 1.032430449278663
$ws.Range("M8").Value = 1.058815588460869
$ws.Range("N8").Value = 1.036412593944449
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027878349007965
$ws.Range("D9").Value = 1.050630611741759
$ws.Range("E9").Value = 1.02795414535692
$ws.Range("F9").Value = 1.054117423217909
$ws.Range("I9").Value = 1.039557466201513
$ws.Range("J9").Value = 1.033728985214864
$ws.Range("K9").Value = 1.053746189778571
$ws.Range("L9").Value = 1.031143195088752
$ws.Range("M9").Value = 1.05722201063851
$ws.Range("N9").Value = 1.035196998944241
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026521861780811
$ws.Range("D10").Value = 1.049562757565078
$ws.Range("E10").Value = 1.026807683348021
$ws.Range("F10").Value = 1.052774390813806
$ws.Range("I10").Value = 1.039232430503465
$ws.Range("J10").Value = 1.032919822660341
$ws.Range("K10").Value = 1.052960587737862
$ws.Range("L10").Value = 1.030286449722943
$ws.Range("M10").Value = 1.056161107920138
$ws.Range("N10").Value = 1.034386687286078
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025935303004988
$ws.Range("D11").Value = 1.049100908049409
$ws.Range("E11").Value = 1.026312430204611
$ws.Range("F11").Value = 1.052194006845945
$ws.Range("I11").Value = 1.03909058237148
$ws.Range("J11").Value = 1.032569480577498
$ws.Range("K11").Value = 1.052620114277299
$ws.Range("L11").Value = 1.029915822684492
$ws.Range("M11").Value = 1.055702096733367
$ws.Range("N11").Value = 1.034035847677293
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025717551843221
$ws.Range("D12").Value = 1.048929439250464
$ws.Range("E12").Value = 1.026128648931125
$ws.Range("F12").Value = 1.051978602224546
$ws.Range("I12").Value = 1.039037728238529
$ws.Range("J12").Value = 1.032439353656103
$ws.Range("K12").Value = 1.052493603168157
$ws.Range("L12").Value = 1.029778209020786
$ws.Range("M12").Value = 1.05553165621906
$ws.Range("N12").Value = 1.033905535960744
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025764254624596
$ws.Range("D13").Value = 1.048966216119988
$ws.Range("E13").Value = 1.026168062566148
$ws.Range("F13").Value = 1.052024799260671
$ws.Range("I13").Value = 1.039049073104293
$ws.Range("J13").Value = 1.032467266045252
$ws.Range("K13").Value = 1.052520742204665
$ws.Range("L13").Value = 1.029807725168017
$ws.Range("M13").Value = 1.055568213697672
$ws.Range("N13").Value = 1.033933487988687
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025917301111582
$ws.Range("D14").Value = 1.049086732688116
$ws.Range("E14").Value = 1.02629723515074
$ws.Range("F14").Value = 1.052176197835186
$ws.Range("I14").Value = 1.039086216803002
$ws.Range("J14").Value = 1.032558724122248
$ws.Range("K14").Value = 1.05260965772147
$ws.Range("L14").Value = 1.029904446393105
$ws.Range("M14").Value = 1.055688006901893
$ws.Range("N14").Value = 1.034025075946643
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026011614426275
$ws.Range("D15").Value = 1.049160997931999
$ws.Range("E15").Value = 1.026376846252097
$ws.Range("F15").Value = 1.052269502848033
$ws.Range("I15").Value = 1.039109080361946
$ws.Range("J15").Value = 1.032615075254515
$ws.Range("K15").Value = 1.05266443570198
$ws.Range("L15").Value = 1.029964046687479
$ws.Range("M15").Value = 1.055761823006545
$ws.Range("N15").Value = 1.034081507103981
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02656080691522
$ws.Range("D16").Value = 1.049593420492272
$ws.Range("E16").Value = 1.026840576511092
$ws.Range("F16").Value = 1.052812933541268
$ws.Range("I16").Value = 1.039241821261923
$ws.Range("J16").Value = 1.032943074452026
$ws.Range("K16").Value = 1.052983177580469
$ws.Range("L16").Value = 1.030311054484706
$ws.Range("M16").Value = 1.056191578842068
$ws.Range("N16").Value = 1.034409972097975
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026905518414845
$ws.Range("D17").Value = 1.049864813144667
$ws.Range("E17").Value = 1.027131777316082
$ws.Range("F17").Value = 1.053154124580674
$ws.Range("I17").Value = 1.039324790613781
$ws.Range("J17").Value = 1.0331488286796
$ws.Range("K17").Value = 1.053183035883356
$ws.Range("L17").Value = 1.03052881771183
$ws.Range("M17").Value = 1.056461252617166
$ws.Range("N17").Value = 1.034616018520152
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027106660782535
$ws.Range("D18").Value = 1.050023163714142
$ws.Range("E18").Value = 1.027301742835197
$ws.Range("F18").Value = 1.053353247139682
$ws.Range("I18").Value = 1.039373078520179
$ws.Range("J18").Value = 1.033268844520641
$ws.Range("K18").Value = 1.053299580563385
$ws.Range("L18").Value = 1.03065586894325
$ws.Range("M18").Value = 1.056618584041577
$ws.Range("N18").Value = 1.034736204797454
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027175258318974
$ws.Range("D19").Value = 1.05007716590252
$ws.Range("E19").Value = 1.027359715804398
$ws.Range("F19").Value = 1.053421161658917
$ws.Range("I19").Value = 1.039389525327621
$ws.Range("J19").Value = 1.033309767275491
$ws.Range("K19").Value = 1.053339314278748
$ws.Range("L19").Value = 1.030699195795785
$ws.Range("M19").Value = 1.056672235926303
$ws.Range("N19").Value = 1.034777185667309
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.026868526057673
$ws.Range("D20").Value = 1.049835689907891
$ws.Range("E20").Value = 1.027100522523927
$ws.Range("F20").Value = 1.053117506454587
$ws.Range("I20").Value = 1.039315899820496
$ws.Range("J20").Value = 1.033126752892092
$ws.Range("K20").Value = 1.053161595998378
$ws.Range("L20").Value = 1.030505450291522
$ws.Range("M20").Value = 1.056432315513263
$ws.Range("N20").Value = 1.034593911382494
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025872229332926
$ws.Range("D21").Value = 1.049051241266009
$ws.Range("E21").Value = 1.026259192122901
$ws.Range("F21").Value = 1.052131609859901
$ws.Range("I21").Value = 1.039075283470532
$ws.Range("J21").Value = 1.03253179183098
$ws.Range("K21").Value = 1.052583475524598
$ws.Range("L21").Value = 1.029875962910172
$ws.Range("M21").Value = 1.055652729211682
$ws.Range("N21").Value = 1.033998105408431
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025246529665518
$ws.Range("D22").Value = 1.048558506486701
$ws.Range("E22").Value = 1.025731243788805
$ws.Range("F22").Value = 1.051512755736616
$ws.Range("I22").Value = 1.038923041681158
$ws.Range("J22").Value = 1.032157749890609
$ws.Range("K22").Value = 1.052219733722823
$ws.Range("L22").Value = 1.029480490630627
$ws.Range("M22").Value = 1.055162901069379
$ws.Range("N22").Value = 1.033623532285602
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025578157002016
$ws.Range("D23").Value = 1.048819668495095
$ws.Range("E23").Value = 1.026011021057933
$ws.Range("F23").Value = 1.051840724881737
$ws.Range("I23").Value = 1.039003838397019
$ws.Range("J23").Value = 1.032356033031137
$ws.Range("K23").Value = 1.052412583806625
$ws.Range("L23").Value = 1.02969010803748
$ws.Range("M23").Value = 1.055422536547739
$ws.Range("N23").Value = 1.033822097010933
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026885241070747
$ws.Range("D24").Value = 1.049848849285081
$ws.Range("E24").Value = 1.027114644870156
$ws.Range("F24").Value = 1.053134052264962
$ws.Range("I24").Value = 1.039319917517084
$ws.Range("J24").Value = 1.033136727981451
$ws.Range("K24").Value = 1.053171283850183
$ws.Range("L24").Value = 1.03051600891917
$ws.Range("M24").Value = 1.056445390836133
$ws.Range("N24").Value = 1.034603900637624
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028405093692672
$ws.Range("D25").Value = 1.051045182893244
$ws.Range("E25").Value = 1.028399766321321
$ws.Range("F25").Value = 1.054639257747924
$ws.Range("I25").Value = 1.039682521181665
$ws.Range("J25").Value = 1.034042789971398
$ws.Range("K25").Value = 1.054050559253962
$ws.Range("L25").Value = 1.031475735506696
$ws.Range("M25").Value = 1.057633734285827
$ws.Range("N25").Value = 1.035511249339523
